$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark that used to sit after "(not bytes)" ---
$d.Bookmarks("_GoBack").Delete()

# --- 2. Insert the four new "function call" / "function call with" doc paragraphs
#         right after the "call stack" section and before "Core Properties" ---
$anchorPara = $d.Paragraphs(28)
$anchorRange = $anchorPara.Range
$anchorRange.InsertParagraphAfter()
$insertPara = $d.Paragraphs(29)
$insertRange = $insertPara.Range
$insertRange.Collapse(1)
$xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>function</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="7030A0"/></w:rPr><w:t xml:space="preserve">call </w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">function </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>param</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Invokes a function variable </w:t></w:r><w:r><w:t>param</w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>function</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="7030A0"/></w:rPr><w:t xml:space="preserve">call </w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">function </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>param</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:color w:val="7030A0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="7030A0"/></w:rPr><w:t>with</w:t></w:r><w:r><w:rPr><w:color w:val="7030A0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>param</w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Invokes a function variable </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>param</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, forwarding all </w:t></w:r><w:r><w:t>params</w:t></w:r><w:r><w:t xml:space="preserve"> to the function.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($xmlPayload)

# --- 3. Re-apply the "CodeChar" character style to the "param"/"params" mentions
#         (InsertXML does not preserve w:rStyle, so patch it up afterwards) ---
$codeParaIndexes = @(30, 33)
$bodyPara = $d.Paragraphs(30)
$bodyRange = $bodyPara.Range
$relStart = $bodyRange.Start + 28
$paramRange = $d.Range($relStart, $relStart + 5)
Write-Output "paraB sub: [$($paramRange.Text)]"
$paramRange.Style = "CodeChar"

$callParaWith = $d.Paragraphs(32)
$callWithRange = $callParaWith.Range
$relStart2 = $callWithRange.Start + 28
$paramRange2 = $d.Range($relStart2, $relStart2 + 5)
Write-Output "paraD sub1: [$($paramRange2.Text)]"
$paramRange2.Style = "CodeChar"

$relStart3 = $callWithRange.Start + 50
$paramsRange = $d.Range($relStart3, $relStart3 + 6)
Write-Output "paraD sub2: [$($paramsRange.Text)]"
$paramsRange.Style = "CodeChar"
